# Commit: "Add files via upload" — re-upload of dir_GRUPOS12.xlsx with a
# corrected category label in cell A1 of sheet "Area".
#
# The shared string "GRUPO12" is fixed to "GRUPOS12" (matches the
# file's own name, dir_GRUPOS12.xlsx). Re-select A1 (the sheet's default
# / top-left cell) so the saved view doesn't keep pointing at the old
# A2 selection left over from the previous save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Area")

$ws.Range("A1").Value = "GRUPOS12"
$ws.Range("A1").Select()
